$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 16,16
$data[0,0] = [double]"7.824700036814042e-23"
$data[0,1] = [double]"-7.90996668856747e-25"
$data[0,2] = [double]"7.909969537698007e-25"
$data[0,3] = [double]"-2.649211851974342e-09"
$data[0,4] = [double]"3.973817138478548e-09"
$data[0,5] = [double]"3.386218234753567e-17"
$data[0,6] = [double]"9.474734366472908e-17"
$data[0,7] = [double]"4.605773699537952e-17"
$data[0,8] = [double]"3.301624516180634e-17"
$data[0,9] = [double]"1.410462627063155e-17"
$data[0,10] = [double]"1"
$data[0,11] = [double]"4.652535058273291e-18"
$data[0,12] = [double]"-4.652532846544745e-18"
$data[0,13] = [double]"-3.300154585485897e-09"
$data[0,14] = [double]"3.313963346709536e-09"
$data[0,15] = [double]"2.483635693623413e-09"
$data[1,0] = [double]"-1.323383709164101e-16"
$data[1,1] = [double]"-7.975649896467908e-18"
$data[1,2] = [double]"7.97565192460229e-18"
$data[1,3] = [double]"1.119538207566623e-16"
$data[1,4] = [double]"-2.423577456099424e-16"
$data[1,5] = [double]"2.838440958027348e-09"
$data[1,6] = [double]"1.51152175241867e-16"
$data[1,7] = [double]"7.574508932552929e-17"
$data[1,8] = [double]"1.049000433672486e-09"
$data[1,9] = [double]"5.440030087187521e-09"
$data[1,10] = [double]"4.675079193501772e-09"
$data[1,11] = [double]"3.177607318161294e-17"
$data[1,12] = [double]"-3.17760589743568e-17"
$data[1,13] = [double]"0.7059034665040979"
$data[1,14] = [double]"-0.7088571622063429"
$data[1,15] = [double]"-1.700767088962572e-17"
$data[2,0] = [double]"3.712191947308202e-18"
$data[2,1] = [double]"-5.693259127314406e-18"
$data[2,2] = [double]"5.693260337500187e-18"
$data[2,3] = [double]"2.720025923924078e-16"
$data[2,4] = [double]"-7.591266567732768e-17"
$data[2,5] = [double]"1.848287136065878e-09"
$data[2,6] = [double]"-2.02693863820576e-17"
$data[2,7] = [double]"-1.986908646456741e-09"
$data[2,8] = [double]"-2.161579654309105e-17"
$data[2,9] = [double]"-9.234323580264316e-18"
$data[2,10] = [double]"2.483635852652983e-09"
$data[2,11] = [double]"-1.873275420810524e-09"
$data[2,12] = [double]"1.873275507322166e-09"
$data[2,13] = [double]"6.356098469985036e-17"
$data[2,14] = [double]"-6.382646608192234e-17"
$data[2,15] = [double]"-1"
$data[3,0] = [double]"1.986909167804789e-09"
$data[3,1] = [double]"-2.809912852870589e-09"
$data[3,2] = [double]"2.80991345529105e-09"
$data[3,3] = [double]"1.015936161921454e-16"
$data[3,4] = [double]"2.654952713007304e-16"
$data[3,5] = [double]"1"
$data[3,6] = [double]"-2.093756012172893e-15"
$data[3,7] = [double]"-1.425667971720675e-15"
$data[3,8] = [double]"3.947404198563188e-17"
$data[3,9] = [double]"1.602277424961643e-16"
$data[3,10] = [double]"-5.172259095606808e-17"
$data[3,11] = [double]"4.822626264909301e-17"
$data[3,12] = [double]"-4.822639054643746e-17"
$data[3,13] = [double]"-2.003665282216909e-09"
$data[3,14] = [double]"2.012049173145372e-09"
$data[3,15] = [double]"1.848287108318087e-09"
$data[4,0] = [double]"8.947778681589694e-18"
$data[4,1] = [double]"1.412579647592062e-22"
$data[4,2] = [double]"-1.412131841458648e-22"
$data[4,3] = [double]"9.940058641573251e-16"
$data[4,4] = [double]"1"
$data[4,5] = [double]"2.914763894239679e-18"
$data[4,6] = [double]"7.947634600841053e-09"
$data[4,7] = [double]"6.623028813236063e-09"
$data[4,8] = [double]"2.630620234079839e-09"
$data[4,9] = [double]"-2.370518208077011e-09"
$data[4,10] = [double]"-3.973817157781597e-09"
$data[4,11] = [double]"-1.654199140042561e-17"
$data[4,12] = [double]"-5.78985487380776e-17"
$data[4,13] = [double]"-2.141862327812799e-16"
$data[4,14] = [double]"1.940335268691096e-16"
$data[4,15] = [double]"-2.776684781839089e-17"
$data[5,0] = [double]"1.655756799150126e-09"
$data[5,1] = [double]"-2.319129235922194e-16"
$data[5,2] = [double]"2.687905912594896e-16"
$data[5,3] = [double]"-3.70386956391966e-25"
$data[5,4] = [double]"4.675079167707451e-09"
$data[5,5] = [double]"-7.002582743014138e-17"
$data[5,6] = [double]"2.024035923506882e-16"
$data[5,7] = [double]"1.703122026984514e-16"
$data[5,8] = [double]"-0.5626899833592589"
$data[5,9] = [double]"0.5070541330848104"
$data[5,10] = [double]"-1.857791010823294e-17"
$data[5,11] = [double]"3.384738009691889e-16"
$data[5,12] = [double]"5.276247773255998e-16"
$data[5,13] = [double]"9.555430008574439e-12"
$data[5,14] = [double]"5.619809096451444e-09"
$data[5,15] = [double]"4.861728531215976e-25"
$data[6,0] = [double]"1.848287189372453e-09"
$data[6,1] = [double]"2.520148960005804e-18"
$data[6,2] = [double]"5.118488226892126e-18"
$data[6,3] = [double]"6.473691970266424e-24"
$data[6,4] = [double]"6.623028820918895e-09"
$data[6,5] = [double]"-2.258308033597706e-16"
$data[6,6] = [double]"-8.509938411708759e-17"
$data[6,7] = [double]"-1"
$data[6,8] = [double]"1.440198129235783e-16"
$data[6,9] = [double]"1.634112224055599e-16"
$data[6,10] = [double]"1.480427496530772e-17"
$data[6,11] = [double]"5.61982610430406e-09"
$data[6,12] = [double]"5.619826513914061e-09"
$data[6,13] = [double]"2.888217771234649e-16"
$data[6,14] = [double]"-2.865149209558429e-16"
$data[6,15] = [double]"1.986908646275665e-09"
$data[7,0] = [double]"1"
$data[7,1] = [double]"2.809913475632295e-09"
$data[7,2] = [double]"2.809912833901234e-09"
$data[7,3] = [double]"3.313932281341595e-23"
$data[7,4] = [double]"-2.895774958993753e-17"
$data[7,5] = [double]"-1.986908568650694e-09"
$data[7,6] = [double]"-2.886728162108907e-17"
$data[7,7] = [double]"1.848287093626292e-09"
$data[7,8] = [double]"9.316779239895282e-10"
$data[7,9] = [double]"-8.395586027018407e-10"
$data[7,10] = [double]"-5.825153784618841e-22"
$data[7,11] = [double]"8.935931456099475e-17"
$data[7,12] = [double]"-1.113175472854856e-16"
$data[7,13] = [double]"-2.176516936295965e-16"
$data[7,14] = [double]"2.076513655140275e-16"
$data[7,15] = [double]"-3.672377633393241e-18"
$data[8,0] = [double]"4.054101869316e-20"
$data[8,1] = [double]"-2.948619879048495e-16"
$data[8,2] = [double]"2.948620914110275e-16"
$data[8,3] = [double]"-1"
$data[8,4] = [double]"4.437151675652741e-16"
$data[8,5] = [double]"-3.019015714223487e-18"
$data[8,6] = [double]"3.973817273692584e-09"
$data[8,7] = [double]"5.990407309556648e-18"
$data[8,8] = [double]"-2.433353792460573e-09"
$data[8,9] = [double]"-2.537103410996985e-09"
$data[8,10] = [double]"-2.649211313002365e-09"
$data[8,11] = [double]"1.756195667610035e-09"
$data[8,12] = [double]"-1.756195726033613e-09"
$data[8,13] = [double]"-9.100017487398644e-17"
$data[8,14] = [double]"7.481045121215513e-17"
$data[8,15] = [double]"-3.801201230001748e-17"
$data[9,0] = [double]"4.120744009782603e-16"
$data[9,1] = [double]"-2.007080659497989e-09"
$data[9,2] = [double]"2.007081081565148e-09"
$data[9,3] = [double]"-2.943568364852842e-09"
$data[9,4] = [double]"3.299371215187659e-24"
$data[9,5] = [double]"-1.306922109869636e-17"
$data[9,6] = [double]"1.7450977035167e-17"
$data[9,7] = [double]"-2.005011761696782e-16"
$data[9,8] = [double]"0.8266680002438446"
$data[9,9] = [double]"0.8619142104186539"
$data[9,10] = [double]"-5.366951929102401e-17"
$data[9,11] = [double]"-3.223745476446737e-16"
$data[9,12] = [double]"5.806604262027855e-16"
$data[9,13] = [double]"9.555298737089613e-12"
$data[9,14] = [double]"5.6198089657278e-09"
$data[9,15] = [double]"1.046948165736032e-24"
$data[10,0] = [double]"2.146244774684155e-16"
$data[10,1] = [double]"-1.702977500901151e-09"
$data[10,2] = [double]"1.70297783893299e-09"
$data[10,3] = [double]"-2.483635807844583e-09"
$data[10,4] = [double]"1.86470416435836e-24"
$data[10,5] = [double]"-1.821701640826346e-16"
$data[10,6] = [double]"1.213106796442075e-16"
$data[10,7] = [double]"-3.743229428101391e-16"
$data[10,8] = [double]"-3.217809688236831e-16"
$data[10,9] = [double]"-2.409134217358464e-16"
$data[10,10] = [double]"-6.579676175932491e-18"
$data[10,11] = [double]"-0.7071067724033828"
$data[10,12] = [double]"0.7071067900728011"
$data[10,13] = [double]"-3.096184749165145e-16"
$data[10,14] = [double]"-3.083263362408814e-16"
$data[10,15] = [double]"2.649211528367558e-09"
$data[11,0] = [double]"4.040313890601668e-17"
$data[11,1] = [double]"-0.7071067110086425"
$data[11,2] = [double]"0.7071068513644451"
$data[11,3] = [double]"4.313343468271358e-16"
$data[11,4] = [double]"-1.946361808112078e-24"
$data[11,5] = [double]"-3.973817292804073e-09"
$data[11,6] = [double]"2.188007680662862e-19"
$data[11,7] = [double]"-5.732538826621778e-19"
$data[11,8] = [double]"-2.34644844317791e-09"
$data[11,9] = [double]"-2.446492450713176e-09"
$data[11,10] = [double]"2.870679841413234e-25"
$data[11,11] = [double]"1.702977630388271e-09"
$data[11,12] = [double]"-1.702977708783019e-09"
$data[11,13] = [double]"-6.269234402664087e-17"
$data[11,14] = [double]"-7.835519445513556e-17"
$data[11,15] = [double]"-5.673726337629256e-18"
$data[12,0] = [double]"6.715470931974654e-16"
$data[12,1] = [double]"-1.017640529324797e-17"
$data[12,2] = [double]"-1.156542670386365e-17"
$data[12,3] = [double]"-3.973817292551323e-09"
$data[12,4] = [double]"7.947634585102668e-09"
$data[12,5] = [double]"-2.299588223035076e-15"
$data[12,6] = [double]"-1"
$data[12,7] = [double]"-8.500700718805016e-16"
$data[12,8] = [double]"1.536952353297894e-16"
$data[12,9] = [double]"-9.025797794251441e-17"
$data[12,10] = [double]"5.263741594400079e-17"
$data[12,11] = [double]"4.68318879815399e-09"
$data[12,12] = [double]"4.683188383694625e-09"
$data[12,13] = [double]"2.08495304995301e-09"
$data[12,14] = [double]"2.076252123338593e-09"
$data[12,15] = [double]"1.523444410732731e-24"
$data[13,0] = [double]"8.69903843991277e-18"
$data[13,1] = [double]"-1.170797254732926e-09"
$data[13,2] = [double]"-1.170797040890147e-09"
$data[13,3] = [double]"-1.16972028700551e-17"
$data[13,4] = [double]"2.711370712964966e-17"
$data[13,5] = [double]"-6.386081660972935e-19"
$data[13,6] = [double]"-2.943568364277768e-09"
$data[13,7] = [double]"-1.203193787887785e-18"
$data[13,8] = [double]"1.049000622657662e-09"
$data[13,9] = [double]"5.440030157405524e-09"
$data[13,10] = [double]"-9.167840513775201e-26"
$data[13,11] = [double]"7.014040105454907e-17"
$data[13,12] = [double]"-2.584020180877375e-18"
$data[13,13] = [double]"-0.7083080516113721"
$data[13,14] = [double]"-0.7053520564858163"
$data[13,15] = [double]"4.640120157007941e-33"
$data[14,0] = [double]"8.730491183095756e-18"
$data[14,1] = [double]"-1.702977870764936e-09"
$data[14,2] = [double]"-1.702977466731791e-09"
$data[14,3] = [double]"-1.603796174762036e-17"
$data[14,4] = [double]"5.263741291526675e-17"
$data[14,5] = [double]"-5.117247495282997e-23"
$data[14,6] = [double]"-6.623028820918865e-09"
$data[14,7] = [double]"-7.947634585102764e-09"
$data[14,8] = [double]"-2.920747394811565e-16"
$data[14,9] = [double]"4.124799415372144e-17"
$data[14,10] = [double]"2.418545036754789e-24"
$data[14,11] = [double]"-0.707106789969712"
$data[14,12] = [double]"-0.7071067723002938"
$data[14,13] = [double]"-1.720297383699401e-17"
$data[14,14] = [double]"-1.713118541548484e-17"
$data[14,15] = [double]"-1.228206301356223e-17"
$data[15,0] = [double]"3.973817292551331e-09"
$data[15,1] = [double]"-0.7071068513644454"
$data[15,2] = [double]"-0.707106711008643"
$data[15,3] = [double]"-3.207658800987072e-23"
$data[15,4] = [double]"-3.668979670910397e-25"
$data[15,5] = [double]"2.36868352585634e-17"
$data[15,6] = [double]"3.620562805373819e-17"
$data[15,7] = [double]"2.107625470069782e-17"
$data[15,8] = [double]"-1.179710779764038e-16"
$data[15,9] = [double]"-1.802916187810417e-16"
$data[15,10] = [double]"-1.891876555947446e-30"
$data[15,11] = [double]"1.702977690718497e-09"
$data[15,12] = [double]"1.702977648163928e-09"
$data[15,13] = [double]"1.172786159977733e-09"
$data[15,14] = [double]"1.167891749749972e-09"
$data[15,15] = [double]"5.326319502923056e-17"
$ws.Range("A2:P17").Value = $data
Write-Output "done"
